$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($addr, $val)
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue 'D2' '34.468.81'
Set-TextValue 'E2' '  +0.37%  '
Set-TextValue 'D3' '1.806.86'
Set-TextValue 'E3' '  +0.11%  '
Set-TextValue 'E4' '  +0.04%  '
Set-TextValue 'D5' '225.60'
Set-TextValue 'E5' '  -0.83%  '
Set-TextValue 'E6' '  +2.28%  '
Set-TextValue 'E7' '  +0.01%  '
Set-TextValue 'D8' '38.20'
Set-TextValue 'E8' '  +5.69%  '
Set-TextValue 'E9' '  -4.44%  '
Set-TextValue 'D10' '0.0671'
Set-TextValue 'E10' '  -3.30%  '
Set-TextValue 'E11' '  +0.82%  '
Set-TextValue 'D12' '2.069.72'
Set-TextValue 'D13' '11.10'
Set-TextValue 'D14' '1.817.22'
Set-TextValue 'E14' '  +0.61%  '
Set-TextValue 'D15' '34.459.27'
Set-TextValue 'E15' '  +0.35%  '
Set-TextValue 'E16' '  -2.64%  '
Set-TextValue 'D17' '4.40'
Set-TextValue 'E17' '  -2.31%  '
Set-TextValue 'D18' '67.98'
Set-TextValue 'E18' '  -1.62%  '
Set-TextValue 'D19' '242.04'
Set-TextValue 'D20' '0.0₃0768'
Set-TextValue 'E20' '  -3.38%  '
Set-TextValue 'D21' '11.10'
Set-TextValue 'E21' '  -4.28%  '
Set-TextValue 'E23' '  -2.04%  '
Set-TextValue 'E24' '  +3.62%  '
Set-TextValue 'D25' '170.49'
Set-TextValue 'E25' '  -0.73%  '
Set-TextValue 'E26' '  -3.60%  '
Set-TextValue 'D27' '17.49'
Set-TextValue 'E27' '  +3.52%  '
Set-TextValue 'E28' '  +1.46%  '
Set-TextValue 'E29' '  +0.00%  '
Set-TextValue 'E30' '  -1.06%  '
Set-TextValue 'E31' '  -2.29%  '
Set-TextValue 'D32' '0.0513'
Set-TextValue 'E32' '  -3.71%  '
Set-TextValue 'E33' '  -5.11%  '
Set-TextValue 'E34' '  -0.82%  '
Set-TextValue 'D35' '1.332.37'
Set-TextValue 'E35' '  -4.21%  '
Set-TextValue 'E36' '  -5.09%  '
Set-TextValue 'D37' '1.05'
Set-TextValue 'E37' '  -1.01%  '
Set-TextValue 'E38' '  -1.13%  '
Set-TextValue 'E39' '  -6.02%  '
Set-TextValue 'E40' '  +1.57%  '
Set-TextValue 'E41' '  -1.58%  '
Set-TextValue 'D42' '81.59'
Set-TextValue 'E42' '  -1.02%  '
Set-TextValue 'D43' '0.946'
Set-TextValue 'E43' '  -2.12%  '
Set-TextValue 'E44' '  -1.13%  '
Set-TextValue 'D45' '13.60'
Set-TextValue 'E45' '  +0.98%  '
Set-TextValue 'D46' '0.0511'
Set-TextValue 'E46' '  +1.64%  '
Set-TextValue 'D47' '1.969.89'
Set-TextValue 'E47' '  +0.17%  '
Set-TextValue 'D48' '5.76'
Set-TextValue 'E48' '  -4.19%  '
Set-TextValue 'E49' '  -0.03%  '
Set-TextValue 'D50' '101.93'
Set-TextValue 'E50' '  -2.46%  '
Set-TextValue 'E51' '  -5.15%  '
